$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "43.738.02"
Set-TextValue "E2" "  +0.17%  "
Set-TextValue "D3" "2.348.80"
Set-TextValue "E3" "  +0.35%  "
Set-TextValue "E4" "  +0.23%  "
Set-TextValue "D5" "0.673"
Set-TextValue "E5" "  +3.33%  "
Set-TextValue "D6" "235.09"
Set-TextValue "E6" "  +0.65%  "
Set-TextValue "D7" "73.39"
Set-TextValue "E7" "  +9.53%  "
Set-TextValue "E8" "  +0.03%  "
Set-TextValue "D9" "0.542"
Set-TextValue "E9" "  +18.30%  "
Set-TextValue "D10" "0.0985"
Set-TextValue "E10" "  +1.26%  "
Set-TextValue "D11" "28.13"
Set-TextValue "E11" "  +3.65%  "
Set-TextValue "E12" "  +1.96%  "
Set-TextValue "D13" "2.701.64"
Set-TextValue "E13" "  +0.45%  "
Set-TextValue "E14" "  +6.83%  "
Set-TextValue "E15" "  +6.98%  "
Set-TextValue "D16" "0.889"
Set-TextValue "E16" "  +4.24%  "
Set-TextValue "D17" "2.323.67"
Set-TextValue "E17" "  -1.07%  "
Set-TextValue "D18" "43.707.43"
Set-TextValue "E18" "  +0.11%  "
Set-TextValue "E19" "  +3.39%  "
Set-TextValue "D20" "77.00"
Set-TextValue "E20" "  +3.55%  "
Set-TextValue "D21" "6.39"
Set-TextValue "E21" "  +1.88%  "
Set-TextValue "D22" "252.92"
Set-TextValue "E22" "  +1.19%  "
Set-TextValue "D24" "3.76"
Set-TextValue "E24" "  -1.43%  "
Set-TextValue "E25" "  +1.60%  "
Set-TextValue "D26" "10.55"
Set-TextValue "E26" "  +5.51%  "
Set-TextValue "D27" "2.30"
Set-TextValue "E27" "  +0.78%  "
Set-TextValue "E28" "  +0.30%  "
Set-TextValue "E29" "  +8.39%  "
Set-TextValue "D30" "172.31"
Set-TextValue "E30" "  -1.56%  "
Set-TextValue "E31" "  -0.03%  "
Set-TextValue "E32" "  +4.69%  "
Set-TextValue "E33" "  +2.60%  "
Set-TextValue "D34" "0.0713"
Set-TextValue "E34" "  +3.31%  "
Set-TextValue "E35" "  +3.31%  "
Set-TextValue "E36" "  +5.59%  "
Set-TextValue "D37" "2.40"
Set-TextValue "E37" "  -3.84%  "
Set-TextValue "E38" "  -3.00%  "
Set-TextValue "E39" "  +5.75%  "
Set-TextValue "D40" "19.31"
Set-TextValue "E40" "  +5.27%  "
Set-TextValue "E41" "  +0.28%  "
Set-TextValue "D42" "8.86"
Set-TextValue "E42" "  -2.14%  "
Set-TextValue "D43" "0.0978"
Set-TextValue "E43" "  +2.39%  "
Set-TextValue "B44" "TrustWalletToken"
Set-TextValue "C44" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D44" "1.23"
Set-TextValue "E44" "  +1.93%  "
Set-TextValue "B45" "ARBITRUM"
Set-TextValue "C45" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D45" "1.16"
Set-TextValue "E45" "  -1.82%  "
Set-TextValue "D46" "0.183"
Set-TextValue "E46" "  +12.29%  "
Set-TextValue "E47" "  +1.81%  "
Set-TextValue "D48" "97.16"
Set-TextValue "E48" "  -2.52%  "
Set-TextValue "D49" "1.431.49"
Set-TextValue "E49" "  -1.09%  "
Set-TextValue "E50" "  +1.42%  "
Set-TextValue "D51" "2.572.87"
Set-TextValue "E51" "  +0.41%  "
